$d = $word.ActiveDocument

$d.Content.Find.Execute("72+18=90", $true, $false, $false, $false, $false, $true, 1, $false, "61-0=61", 2) | Out-Null
$d.Content.Find.Execute("43+6=49", $true, $false, $false, $false, $false, $true, 1, $false, "86-5=81", 2) | Out-Null
$d.Content.Find.Execute("61-20=41", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=80", 2) | Out-Null
$d.Content.Find.Execute("78+15=93", $true, $false, $false, $false, $false, $true, 1, $false, "52-36=16", 2) | Out-Null
$d.Content.Find.Execute("48+29=77", $true, $false, $false, $false, $false, $true, 1, $false, "48-29=19", 2) | Out-Null
$d.Content.Find.Execute("18+50=68", $true, $false, $false, $false, $false, $true, 1, $false, "14-0=14", 2) | Out-Null
$d.Content.Find.Execute("90-0=90", $true, $false, $false, $false, $false, $true, 1, $false, "28+11=39", 2) | Out-Null
$d.Content.Find.Execute("28+28=56", $true, $false, $false, $false, $false, $true, 1, $false, "5+61=66", 2) | Out-Null
$d.Content.Find.Execute("38-11=27", $true, $false, $false, $false, $false, $true, 1, $false, "65-61=4", 2) | Out-Null
$d.Content.Find.Execute("11+71=82", $true, $false, $false, $false, $false, $true, 1, $false, "52+22=74", 2) | Out-Null
$d.Content.Find.Execute("89-30=59", $true, $false, $false, $false, $false, $true, 1, $false, "8+56=64", 2) | Out-Null
$d.Content.Find.Execute("60-27=33", $true, $false, $false, $false, $false, $true, 1, $false, "35+21=56", 2) | Out-Null
$d.Content.Find.Execute("48-40=8", $true, $false, $false, $false, $false, $true, 1, $false, "15+70=85", 2) | Out-Null
$d.Content.Find.Execute("49+33=82", $true, $false, $false, $false, $false, $true, 1, $false, "4+19=23", 2) | Out-Null
$d.Content.Find.Execute("80-31=49", $true, $false, $false, $false, $false, $true, 1, $false, "75+21=96", 2) | Out-Null
$d.Content.Find.Execute("28+67=95", $true, $false, $false, $false, $false, $true, 1, $false, "72-66=6", 2) | Out-Null
$d.Content.Find.Execute("0+74=74", $true, $false, $false, $false, $false, $true, 1, $false, "90-54=36", 2) | Out-Null
$d.Content.Find.Execute("49-44=5", $true, $false, $false, $false, $false, $true, 1, $false, "90-39=51", 2) | Out-Null
$d.Content.Find.Execute("51-18=33", $true, $false, $false, $false, $false, $true, 1, $false, "72+12=84", 2) | Out-Null
$d.Content.Find.Execute("74+1=75", $true, $false, $false, $false, $false, $true, 1, $false, "14+3=17", 2) | Out-Null
$d.Content.Find.Execute("10+79=89", $true, $false, $false, $false, $false, $true, 1, $false, "33+64=97", 2) | Out-Null
$d.Content.Find.Execute("43+4=47", $true, $false, $false, $false, $false, $true, 1, $false, "92-87=5", 2) | Out-Null
$d.Content.Find.Execute("91-34=57", $true, $false, $false, $false, $false, $true, 1, $false, "96-91=5", 2) | Out-Null
$d.Content.Find.Execute("71-25=46", $true, $false, $false, $false, $false, $true, 1, $false, "22+57=79", 2) | Out-Null
$d.Content.Find.Execute("35-8=27", $true, $false, $false, $false, $false, $true, 1, $false, "51-0=51", 2) | Out-Null
$d.Content.Find.Execute("1+84=85", $true, $false, $false, $false, $false, $true, 1, $false, "5+45=50", 2) | Out-Null
$d.Content.Find.Execute("56-1=55", $true, $false, $false, $false, $false, $true, 1, $false, "49-47=2", 2) | Out-Null
$d.Content.Find.Execute("39+0=39", $true, $false, $false, $false, $false, $true, 1, $false, "59+13=72", 2) | Out-Null
$d.Content.Find.Execute("99-69=30", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=81", 2) | Out-Null
$d.Content.Find.Execute("1+90=91", $true, $false, $false, $false, $false, $true, 1, $false, "44-10=34", 2) | Out-Null
$d.Content.Find.Execute("16+9=25", $true, $false, $false, $false, $false, $true, 1, $false, "36-10=26", 2) | Out-Null
$d.Content.Find.Execute("5+83=88", $true, $false, $false, $false, $false, $true, 1, $false, "72-42=30", 2) | Out-Null
$d.Content.Find.Execute("66-49=17", $true, $false, $false, $false, $false, $true, 1, $false, "53-16=37", 2) | Out-Null
$d.Content.Find.Execute("14+27=41", $true, $false, $false, $false, $false, $true, 1, $false, "16+69=85", 2) | Out-Null
$d.Content.Find.Execute("54+45=99", $true, $false, $false, $false, $false, $true, 1, $false, "38+13=51", 2) | Out-Null
$d.Content.Find.Execute("11+26=37", $true, $false, $false, $false, $false, $true, 1, $false, "2+52=54", 2) | Out-Null
$d.Content.Find.Execute("66+18=84", $true, $false, $false, $false, $false, $true, 1, $false, "18-2=16", 2) | Out-Null
$d.Content.Find.Execute("99-74=25", $true, $false, $false, $false, $false, $true, 1, $false, "68+25=93", 2) | Out-Null
$d.Content.Find.Execute("30+28=58", $true, $false, $false, $false, $false, $true, 1, $false, "44+43=87", 2) | Out-Null
$d.Content.Find.Execute("14+40=54", $true, $false, $false, $false, $false, $true, 1, $false, "43+49=92", 2) | Out-Null
$d.Content.Find.Execute("82-59=23", $true, $false, $false, $false, $false, $true, 1, $false, "88-37=51", 2) | Out-Null
$d.Content.Find.Execute("92-12=80", $true, $false, $false, $false, $false, $true, 1, $false, "89-5=84", 2) | Out-Null
$d.Content.Find.Execute("51+22=73", $true, $false, $false, $false, $false, $true, 1, $false, "0+84=84", 2) | Out-Null
$d.Content.Find.Execute("66-5=61", $true, $false, $false, $false, $false, $true, 1, $false, "38+60=98", 2) | Out-Null
$d.Content.Find.Execute("31+49=80", $true, $false, $false, $false, $false, $true, 1, $false, "68+19=87", 2) | Out-Null
$d.Content.Find.Execute("11-3=8", $true, $false, $false, $false, $false, $true, 1, $false, "4+82=86", 2) | Out-Null
$d.Content.Find.Execute("16+2=18", $true, $false, $false, $false, $false, $true, 1, $false, "86-12=74", 2) | Out-Null
$d.Content.Find.Execute("5+4=9", $true, $false, $false, $false, $false, $true, 1, $false, "83-11=72", 2) | Out-Null
$d.Content.Find.Execute("73-33=40", $true, $false, $false, $false, $false, $true, 1, $false, "71-36=35", 2) | Out-Null
$d.Content.Find.Execute("53-25=28", $true, $false, $false, $false, $false, $true, 1, $false, "94-29=65", 2) | Out-Null
$d.Content.Find.Execute("88-80=8", $true, $false, $false, $false, $false, $true, 1, $false, "97-19=78", 2) | Out-Null
$d.Content.Find.Execute("74+5=79", $true, $false, $false, $false, $false, $true, 1, $false, "26-1=25", 2) | Out-Null
$d.Content.Find.Execute("13+62=75", $true, $false, $false, $false, $false, $true, 1, $false, "18+38=56", 2) | Out-Null
$d.Content.Find.Execute("90-55=35", $true, $false, $false, $false, $false, $true, 1, $false, "71+8=79", 2) | Out-Null
$d.Content.Find.Execute("17+19=36", $true, $false, $false, $false, $false, $true, 1, $false, "77-66=11", 2) | Out-Null
$d.Content.Find.Execute("2+41=43", $true, $false, $false, $false, $false, $true, 1, $false, "49-42=7", 2) | Out-Null
$d.Content.Find.Execute("52-31=21", $true, $false, $false, $false, $false, $true, 1, $false, "78-74=4", 2) | Out-Null
$d.Content.Find.Execute("87-35=52", $true, $false, $false, $false, $false, $true, 1, $false, "10+7=17", 2) | Out-Null
$d.Content.Find.Execute("56-55=1", $true, $false, $false, $false, $false, $true, 1, $false, "79-69=10", 2) | Out-Null
$d.Content.Find.Execute("44-19=25", $true, $false, $false, $false, $false, $true, 1, $false, "5+58=63", 2) | Out-Null
$d.Content.Find.Execute("23+42=65", $true, $false, $false, $false, $false, $true, 1, $false, "63+19=82", 2) | Out-Null
$d.Content.Find.Execute("83-53=30", $true, $false, $false, $false, $false, $true, 1, $false, "18+18=36", 2) | Out-Null
$d.Content.Find.Execute("8+80=88", $true, $false, $false, $false, $false, $true, 1, $false, "75-49=26", 2) | Out-Null
$d.Content.Find.Execute("56+4=60", $true, $false, $false, $false, $false, $true, 1, $false, "68+0=68", 2) | Out-Null
$d.Content.Find.Execute("25+68=93", $true, $false, $false, $false, $false, $true, 1, $false, "27+14=41", 2) | Out-Null
$d.Content.Find.Execute("15+34=49", $true, $false, $false, $false, $false, $true, 1, $false, "70-22=48", 2) | Out-Null
$d.Content.Find.Execute("64+29=93", $true, $false, $false, $false, $false, $true, 1, $false, "78-51=27", 2) | Out-Null
$d.Content.Find.Execute("59-23=36", $true, $false, $false, $false, $false, $true, 1, $false, "47-22=25", 2) | Out-Null
$d.Content.Find.Execute("41+25=66", $true, $false, $false, $false, $false, $true, 1, $false, "81-20=61", 2) | Out-Null
$d.Content.Find.Execute("97-18=79", $true, $false, $false, $false, $false, $true, 1, $false, "59+40=99", 2) | Out-Null
$d.Content.Find.Execute("86-31=55", $true, $false, $false, $false, $false, $true, 1, $false, "24+24=48", 2) | Out-Null
$d.Content.Find.Execute("46+20=66", $true, $false, $false, $false, $false, $true, 1, $false, "73-21=52", 2) | Out-Null
$d.Content.Find.Execute("1+93=94", $true, $false, $false, $false, $false, $true, 1, $false, "59+13=72", 2) | Out-Null
$d.Content.Find.Execute("55-38=17", $true, $false, $false, $false, $false, $true, 1, $false, "94-7=87", 2) | Out-Null
$d.Content.Find.Execute("92-82=10", $true, $false, $false, $false, $false, $true, 1, $false, "45+48=93", 2) | Out-Null
$d.Content.Find.Execute("58-45=13", $true, $false, $false, $false, $false, $true, 1, $false, "33-31=2", 2) | Out-Null
$d.Content.Find.Execute("50-46=4", $true, $false, $false, $false, $false, $true, 1, $false, "83-59=24", 2) | Out-Null
$d.Content.Find.Execute("78+12=90", $true, $false, $false, $false, $false, $true, 1, $false, "52+34=86", 2) | Out-Null
$d.Content.Find.Execute("86-76=10", $true, $false, $false, $false, $false, $true, 1, $false, "5+48=53", 2) | Out-Null
$d.Content.Find.Execute("62+16=78", $true, $false, $false, $false, $false, $true, 1, $false, "87-58=29", 2) | Out-Null
$d.Content.Find.Execute("14+41=55", $true, $false, $false, $false, $false, $true, 1, $false, "71+5=76", 2) | Out-Null
$d.Content.Find.Execute("56-0=56", $true, $false, $false, $false, $false, $true, 1, $false, "39-35=4", 2) | Out-Null
$d.Content.Find.Execute("7+71=78", $true, $false, $false, $false, $false, $true, 1, $false, "45+3=48", 2) | Out-Null
$d.Content.Find.Execute("24+43=67", $true, $false, $false, $false, $false, $true, 1, $false, "91-8=83", 2) | Out-Null
$d.Content.Find.Execute("34-14=20", $true, $false, $false, $false, $false, $true, 1, $false, "99-58=41", 2) | Out-Null
$d.Content.Find.Execute("30+42=72", $true, $false, $false, $false, $false, $true, 1, $false, "78-20=58", 2) | Out-Null
$d.Content.Find.Execute("79+5=84", $true, $false, $false, $false, $false, $true, 1, $false, "78-26=52", 2) | Out-Null
$d.Content.Find.Execute("32+43=75", $true, $false, $false, $false, $false, $true, 1, $false, "50+3=53", 2) | Out-Null
$d.Content.Find.Execute("44-17=27", $true, $false, $false, $false, $false, $true, 1, $false, "2+22=24", 2) | Out-Null
$d.Content.Find.Execute("41+35=76", $true, $false, $false, $false, $false, $true, 1, $false, "75-28=47", 2) | Out-Null
$d.Content.Find.Execute("16+59=75", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=6", 2) | Out-Null
$d.Content.Find.Execute("14+19=33", $true, $false, $false, $false, $false, $true, 1, $false, "78-35=43", 2) | Out-Null
$d.Content.Find.Execute("74-15=59", $true, $false, $false, $false, $false, $true, 1, $false, "23+54=77", 2) | Out-Null
$d.Content.Find.Execute("81-49=32", $true, $false, $false, $false, $false, $true, 1, $false, "37+23=60", 2) | Out-Null
$d.Content.Find.Execute("68-56=12", $true, $false, $false, $false, $false, $true, 1, $false, "6+61=67", 2) | Out-Null
$d.Content.Find.Execute("24+75=99", $true, $false, $false, $false, $false, $true, 1, $false, "98-59=39", 2) | Out-Null
$d.Content.Find.Execute("72-8=64", $true, $false, $false, $false, $false, $true, 1, $false, "38-7=31", 2) | Out-Null
$d.Content.Find.Execute("60-52=8", $true, $false, $false, $false, $false, $true, 1, $false, "82-67=15", 2) | Out-Null
$d.Content.Find.Execute("26+59=85", $true, $false, $false, $false, $false, $true, 1, $false, "53+6=59", 2) | Out-Null
$d.Content.Find.Execute("10+19=29", $true, $false, $false, $false, $false, $true, 1, $false, "53+25=78", 2) | Out-Null
